# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Ajo" right before the existing
# row that used to be row 237 (date 2022-06-13 / serial 44685),
# pushing all subsequent rows down by 3 (old 237..266 -> new 240..269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 237:239 - existing rows 237..266 shift to 240..269.
$ws.Rows("237:239").Insert()

# Row 237 (new)
$ws.Range("A237").Value = 9
$ws.Range("B237").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C237").Value = "Metropolitana"
$ws.Range("D237").Value = 44816
$ws.Range("D237").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E237").Value = 13
$ws.Range("F237").Value = 100112003
$ws.Range("G237").Value = "Ajo"
$ws.Range("H237").Value = "Chino"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 430
$ws.Range("K237").Value = 18000
$ws.Range("L237").Value = 18000
$ws.Range("M237").Value = 18000
$ws.Range("N237").Value = "$/caja 10 kilos"
$ws.Range("O237").Value = "China"
$ws.Range("P237").Value = 1800
$ws.Range("Q237").Value = 10
$ws.Range("R237").Value = "Hortaliza"

# Row 238 (new)
$ws.Range("A238").Value = 9
$ws.Range("B238").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C238").Value = "Metropolitana"
$ws.Range("D238").Value = 44816
$ws.Range("D238").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E238").Value = 13
$ws.Range("F238").Value = 100112003
$ws.Range("G238").Value = "Ajo"
$ws.Range("H238").Value = "Chino"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 340
$ws.Range("K238").Value = 22000
$ws.Range("L238").Value = 22000
$ws.Range("M238").Value = 22000
$ws.Range("N238").Value = "$/malla 10 kilos"
$ws.Range("O238").Value = "China"
$ws.Range("P238").Value = 2200
$ws.Range("Q238").Value = 10
$ws.Range("R238").Value = "Hortaliza"

# Row 239 (new)
$ws.Range("A239").Value = 9
$ws.Range("B239").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C239").Value = "Metropolitana"
$ws.Range("D239").Value = 44816
$ws.Range("D239").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E239").Value = 13
$ws.Range("F239").Value = 100112003
$ws.Range("G239").Value = "Ajo"
$ws.Range("H239").Value = "Chino"
$ws.Range("I239").Value = "Segunda"
$ws.Range("J239").Value = 250
$ws.Range("K239").Value = 12000
$ws.Range("L239").Value = 12000
$ws.Range("M239").Value = 12000
$ws.Range("N239").Value = "$/caja 10 kilos"
$ws.Range("O239").Value = "China"
$ws.Range("P239").Value = 1200
$ws.Range("Q239").Value = 10
$ws.Range("R239").Value = "Hortaliza"
